# Added column labels for cartesian/polar matrices
$wb = $excel.ActiveWorkbook

$polarValues = @(
    "user ID",
    "doc ID",
    "inter-stroke time",
    "stroke duration",
    "mid-stroke pressure",
    "mid-stroke area covered",
    "rho start",
    "theta start",
    "rho end",
    "theta end",
    "20 percent drho/dt",
    "50 percent drho/dt",
    "80 percent drho/dt",
    "20 percent dtheta/dt",
    "50 percent dtheta/dt",
    "80 percent dtheta/dt",
    "20 percent d^2rho/dt^2",
    "50 percent d^2rho/dt^2",
    "80 percent d^2rho/dt^2",
    "20 percent d^2theta/dt^2",
    "50 percent d^2theta/dt^2",
    "80 percent d^2theta/dt^2",
    "median drho/dt at last 3 point",
    "median dtheta/dt at last 3 point",
    "rho-displacement",
    "theta-displacement"
)

$cartesianValues = @(
    "user ID",
    "doc ID",
    "inter-stroke time",
    "stroke duration",
    "start x",
    "start y",
    "stop x",
    "stop y",
    "direct end-to-end distance",
    "mean resultant length",
    "up/down/left/right flag",
    "direction of end-to-end line",
    "20 percent pairwise velocity",
    "50 percent pairwaise velocity",
    "80 percent pairwise veolcity",
    "20 percent pairwise acceleration",
    "50 percent pairwise acceleration",
    "80 percent pairwise acceleration",
    "median veolcity at last 3 point",
    "largest deviation from end-to-end line",
    "20 percent deviation from end-to-end line",
    "50 percent deviation from end-to-end line",
    "80 percent deviation from end-to-end line",
    "average direction",
    "length of trajectory",
    "ratio end-to-end distance and length of trajectory",
    "average velocity",
    "median acceleration at first 5 points",
    "mid-stroke pressure",
    "mid-stroke area covered",
    "x-displacment",
    "y-displacement"
)

# Add the two new sheets at the end, in order: polarMatrix, then cartesianMatrix
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPolar = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsPolar.Name = "polarMatrix"

$wsCartesian = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsPolar)
$wsCartesian.Name = "cartesianMatrix"

# Reuse the label-column formatting already used on the other sheets (style with
# the Times New Roman font) by copying the format of an existing label cell.
$wsFormatSource = $wb.Worksheets.Item(1).Range("A1")
$wsFormatSource.Copy()
$wsPolar.Range("A1:A26").PasteSpecial(-4122)
$wsCartesian.Range("A1:A32").PasteSpecial(-4122)

# Fill the polarMatrix sheet (column A)
for ($i = 0; $i -lt $polarValues.Length; $i++) {
    $wsPolar.Cells.Item($i + 1, 1).Value = $polarValues[$i]
}

# Fill the cartesianMatrix sheet (column A)
for ($i = 0; $i -lt $cartesianValues.Length; $i++) {
    $wsCartesian.Cells.Item($i + 1, 1).Value = $cartesianValues[$i]
}

# Set selection on FullLabels to match target state (activeCell/sqref = A55)
$wsFull = $wb.Worksheets.Item("FullLabels")
$wsFull.Select()
$wsFull.Range("A55").Select()

# withoutExtraneous keeps its original selection (E40) in the target diff; only its
# scroll position (topLeftCell) changes, which is a view-only property.
$wsWithout = $wb.Worksheets.Item("withoutExtraneous")
$wsWithout.Select()
$wsWithout.Range("E40").Select()

# Selection / active cell on the new sheets
$wsCartesian.Select()
$wsCartesian.Range("B49").Select()

$wsPolar.Select()
$wsPolar.Range("D29").Select()

# Make polarMatrix the active (selected) tab, matching activeTab="3" (0-indexed 4th sheet)
$wsPolar.Select()
